$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 15.88593828126073
$arr[1,0] = 15.71933865745114
$arr[2,0] = 15.62067092474922
$arr[3,0] = 15.58141980892591
$arr[4,0] = 15.57496125952334
$arr[5,0] = 15.62013763829533
$arr[6,0] = 15.82777007999167
$arr[7,0] = 16.26162181127317
$arr[8,0] = 16.59375540296792
$arr[9,0] = 16.74712541537264
$arr[10,0] = 16.80547722270944
$arr[11,0] = 16.79289879298771
$arr[12,0] = 16.75192087296036
$arr[13,0] = 16.72685478415222
$arr[14,0] = 16.58377358232999
$arr[15,0] = 16.49654356230753
$arr[16,0] = 16.44658907734328
$arr[17,0] = 16.42971442352396
$arr[18,0] = 16.50580718563777
$arr[19,0] = 16.76395007725973
$arr[20,0] = 16.93423309833355
$arr[21,0] = 16.84322382087973
$arr[22,0] = 16.50161848934037
$arr[23,0] = 16.14169025981191
$ws.Range("B2:B25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 8.035776482791084
$arr[1,0] = 7.927829887839084
$arr[2,0] = 7.863331400068546
$arr[3,0] = 7.837529261407999
$arr[4,0] = 7.833274880783852
$arr[5,0] = 7.862981430095066
$arr[6,0] = 7.998207293633159
$arr[7,0] = 8.275997030532865
$arr[8,0] = 8.485705801809679
$arr[9,0] = 8.58188566306633
$arr[10,0] = 8.61838255868115
$arr[11,0] = 8.610519512094113
$arr[12,0] = 8.58488697637776
$arr[13,0] = 8.569195082199602
$arr[14,0] = 8.479432799656893
$arr[15,0] = 8.424540497170385
$arr[16,0] = 8.393043696824524
$arr[17,0] = 8.382393573368278
$arr[18,0] = 8.430376300462617
$arr[19,0] = 8.59241410958438
$arr[20,0] = 8.698739736960595
$arr[21,0] = 8.641964864430165
$arr[22,0] = 8.427737739375926
$arr[23,0] = 8.199710456815344
$ws.Range("C2:C25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 6.293144590183444
$arr[1,0] = 6.295158073387191
$arr[2,0] = 6.297646318034918
$arr[3,0] = 6.298975896282048
$arr[4,0] = 6.299215753525357
$arr[5,0] = 6.297662970397466
$arr[6,0] = 6.293579501650671
$arr[7,0] = 6.295466352744531
$arr[8,0] = 6.302824902069147
$arr[9,0] = 6.307453892168306
$arr[10,0] = 6.309389518068559
$arr[11,0] = 6.308964545775497
$arr[12,0] = 6.307609484541815
$arr[13,0] = 6.30680321842312
$arr[14,0] = 6.302548044707172
$arr[15,0] = 6.300264804167398
$arr[16,0] = 6.299072271966903
$arr[17,0] = 6.298689281307493
$arr[18,0] = 6.300495375048139
$arr[19,0] = 6.308002552778909
$arr[20,0] = 6.313973186922057
$arr[21,0] = 6.310689716821102
$arr[22,0] = 6.300390759648167
$arr[23,0] = 6.293902260774038
$ws.Range("D2:D25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 48.89073580786976
$arr[1,0] = 48.79307069661363
$arr[2,0] = 48.74233128416975
$arr[3,0] = 48.72398406216259
$arr[4,0] = 48.72107844738241
$arr[5,0] = 48.74207440446903
$arr[6,0] = 48.85515181435447
$arr[7,0] = 49.14962991454534
$arr[8,0] = 49.40956775298345
$arr[9,0] = 49.53708441001289
$arr[10,0] = 49.5866845240031
$arr[11,0] = 49.57594423361449
$arr[12,0] = 49.54113886489442
$arr[13,0] = 49.51998983174548
$arr[14,0] = 49.40141904512085
$arr[15,0] = 49.33103977284912
$arr[16,0] = 49.29143333026904
$arr[17,0] = 49.27817395364877
$arr[18,0] = 49.33844147010628
$arr[19,0] = 49.55132660262442
$arr[20,0] = 49.69809775016184
$arr[21,0] = 49.61907158804766
$arr[22,0] = 49.33509249627326
$arr[23,0] = 49.06224952657075
$ws.Range("F2:F25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 3.745829902863885
$arr[1,0] = 3.748934666899126
$arr[2,0] = 3.75094062349457
$arr[3,0] = 3.751783206082527
$arr[4,0] = 3.751924637227957
$arr[5,0] = 3.750951884953973
$arr[6,0] = 3.746879803265869
$arr[7,0] = 3.73968083505161
$arr[8,0] = 3.734865464601659
$arr[9,0] = 3.732776478405726
$arr[10,0] = 3.731999943881974
$arr[11,0] = 3.732166540065267
$arr[12,0] = 3.732712301933506
$arr[13,0] = 3.7330484851946
$arr[14,0] = 3.735004021054275
$arr[15,0] = 3.736229628744927
$arr[16,0] = 3.736944129203631
$arr[17,0] = 3.737187691846358
$arr[18,0] = 3.736098171555065
$arr[19,0] = 3.732551605059285
$arr[20,0] = 3.730318309479981
$arr[21,0] = 3.731502548795066
$arr[22,0] = 3.736157572561923
$arr[23,0] = 3.741544749344564
$ws.Range("G2:G25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 10.99696671043911
$arr[1,0] = 11.01271590913506
$arr[2,0] = 11.02383350078165
$arr[3,0] = 11.02872805318886
$arr[4,0] = 11.02956277810043
$arr[5,0] = 11.02389803638771
$arr[6,0] = 11.00209664182807
$arr[7,0] = 10.97082683512008
$arr[8,0] = 10.95484863694207
$arr[9,0] = 10.94909716926559
$arr[10,0] = 10.94713717162465
$arr[11,0] = 10.94754960205158
$arr[12,0] = 10.94893155216091
$arr[13,0] = 10.94980641448746
$arr[14,0] = 10.95525502078019
$arr[15,0] = 10.958986001606
$arr[16,0] = 10.96127477388246
$arr[17,0] = 10.96207424604163
$arr[18,0] = 10.95857405327209
$arr[19,0] = 10.9485197263847
$arr[20,0] = 10.94321894746354
$arr[21,0] = 10.94593191679182
$arr[22,0] = 10.95875984724096
$arr[23,0] = 10.97805714725441
$ws.Range("J2:J25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 15.4467833546311
$arr[1,0] = 15.3460253989244
$arr[2,0] = 15.28784635449371
$arr[3,0] = 15.2650861293808
$arr[4,0] = 15.26136466625799
$arr[5,0] = 15.28753553620107
$arr[6,0] = 15.41129000680727
$arr[7,0] = 15.68222033666285
$arr[8,0] = 15.89704989844945
$arr[9,0] = 15.99787441061739
$arr[10,0] = 16.03647000811814
$arr[11,0] = 16.02813975750887
$arr[12,0] = 16.0010415392272
$arr[13,0] = 15.98449631500472
$arr[14,0] = 15.89052057252655
$arr[15,0] = 15.83364183088077
$arr[16,0] = 15.80122007355564
$arr[17,0] = 15.79029391109555
$arr[18,0] = 15.83966651926897
$arr[19,0] = 16.00898990673926
$arr[20,0] = 16.12206030098387
$arr[21,0] = 16.06150223938993
$arr[22,0] = 15.83694188761484
$arr[23,0] = 15.60604625192526
$ws.Range("K2:K25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 18.66463188785841
$arr[1,0] = 18.63336699307951
$arr[2,0] = 18.61810819744755
$arr[3,0] = 18.61288486248203
$arr[4,0] = 18.61207772836112
$arr[5,0] = 18.61803372067688
$arr[6,0] = 18.65303736324501
$arr[7,0] = 18.75269430364468
$arr[8,0] = 18.84446250565125
$arr[9,0] = 18.89014602547696
$arr[10,0] = 18.90800236219841
$arr[11,0] = 18.90413206398067
$arr[12,0] = 18.89160396998414
$arr[13,0] = 18.88400239851893
$arr[14,0] = 18.84155544234799
$arr[15,0] = 18.81651779875572
$arr[16,0] = 18.8024880457284
$arr[17,0] = 18.79780184270582
$arr[18,0] = 18.81914474179077
$arr[19,0] = 18.89526873603093
$arr[20,0] = 18.94826186346203
$arr[21,0] = 18.91968506476254
$arr[22,0] = 18.8179559640924
$arr[23,0] = 18.72244784055253
$ws.Range("M2:M25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 24.28350211668339
$arr[1,0] = 24.31690165193041
$arr[2,0] = 24.33925566133304
$arr[3,0] = 24.34882933832147
$arr[4,0] = 24.3504470760597
$arr[5,0] = 24.3393828957975
$arr[6,0] = 24.29463496513944
$arr[7,0] = 24.22154080287013
$arr[8,0] = 24.17678070723196
$arr[9,0] = 24.15836159216114
$arr[10,0] = 24.15166624049944
$arr[11,0] = 24.15309576655874
$arr[12,0] = 24.15780515721607
$arr[13,0] = 24.16072621009837
$arr[14,0] = 24.17802354794222
$arr[15,0] = 24.18913259079474
$arr[16,0] = 24.195705039064
$arr[17,0] = 24.19796174936142
$arr[18,0] = 24.18793109134337
$arr[19,0] = 24.15641430570588
$arr[20,0] = 24.13744591062249
$arr[21,0] = 24.14742050722314
$arr[22,0] = 24.18847371075571
$arr[23,0] = 24.23974490200825
$ws.Range("N2:N25").Value = $arr

Write-Output "Done updating loading_percent values"